$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are treated as text so values like
# "1.00", "0.637", "15.40" keep their exact textual representation
# instead of being auto-converted to numbers by Excel.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D15", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($cellAddr in $priceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Apply updated crypto values scraped by the GitHub Actions job.
$ws.Range("D2").Value = "45.267.99"
$ws.Range("E2").Value = "  -0.82%  "
$ws.Range("D3").Value = "2.366.25"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("D5").Value = "330.33"
$ws.Range("E5").Value = "  +3.84%  "
$ws.Range("D6").Value = "107.28"
$ws.Range("E6").Value = "  -7.24%  "
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "0.615"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "40.99"
$ws.Range("E10").Value = "  -4.31%  "
$ws.Range("D11").Value = "0.0919"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").Value = "8.48"
$ws.Range("E12").Value = "  -2.74%  "
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").Value = "2.726.09"
$ws.Range("E15").Value = "  -1.41%  "
$ws.Range("E16").Value = "  -3.95%  "
$ws.Range("D17").Value = "2.357.90"
$ws.Range("E17").Value = "  -1.96%  "
$ws.Range("D18").Value = "45.233.33"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "15.40"
$ws.Range("E19").Value = "  +13.24%  "
$ws.Range("D20").Value = "7.30"
$ws.Range("E20").Value = "  -3.13%  "
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").Value = "3.66"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("D23").Value = "73.11"
$ws.Range("E23").Value = "  -2.75%  "
$ws.Range("D24").Value = "260.14"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -3.29%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "11.35"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "7.44"
$ws.Range("E28").Value = "  -3.05%  "
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  -2.23%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.0965"
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "22.33"
$ws.Range("E31").Value = "  -2.28%  "
$ws.Range("D32").Value = "36.99"
$ws.Range("E32").Value = "  -8.46%  "
$ws.Range("D33").Value = "167.31"
$ws.Range("E33").Value = "  -3.44%  "
$ws.Range("D34").Value = "2.82"
$ws.Range("E34").Value = "  -4.33%  "
$ws.Range("E35").Value = "  -0.80%  "
$ws.Range("D36").Value = "3.28"
$ws.Range("E36").Value = "  +4.75%  "
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("D38").Value = "4.74"
$ws.Range("E38").Value = "  -5.28%  "
$ws.Range("E39").Value = "  +10.37%  "
$ws.Range("D40").Value = "3.99"
$ws.Range("E40").Value = "  -6.03%  "
$ws.Range("E41").Value = "  -3.38%  "
$ws.Range("D42").Value = "97.20"
$ws.Range("E42").Value = "  -2.70%  "
$ws.Range("D43").Value = "70.22"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "1.889.30"
$ws.Range("E44").Value = "  +14.15%  "
$ws.Range("D45").Value = "0.230"
$ws.Range("E45").Value = "  -5.38%  "
$ws.Range("D46").Value = "6.11"
$ws.Range("E46").Value = "  +4.73%  "
$ws.Range("D47").Value = "12.92"
$ws.Range("E47").Value = "  -6.74%  "
$ws.Range("B48").Value = "ordi"
$ws.Range("C48").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D48").Value = "86.07"
$ws.Range("E48").Value = "  -2.06%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("D50").Value = "112.62"
$ws.Range("E50").Value = "  -3.51%  "
$ws.Range("D51").Value = "9.31"
$ws.Range("E51").Value = "  -2.12%  "
